$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.844.21'
$ws.Range('E2').Value = '  -3.01%  '
$ws.Range('D3').Value = '1.799.01'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.93'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5355'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3870'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07452'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.49'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.089'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.000'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.231'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.521'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.39'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').Value = '1.797.26'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '88.27'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.45%  '
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06513'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.35'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.957'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('D23').Value = '27.876.44'
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.12'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.094'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '156.69'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.29'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('D28').Value = '2.001.95'
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.347'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '121.46'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.118'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1095'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.35%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.653'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.530'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07045'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +9.05%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.2194'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.36%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02277'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.084'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.497'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.51%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '11.32'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6118'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.21%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.165'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('B43').Value = 'WEMIXTOKEN'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.411'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.30'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.670'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5716'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.77%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '125.01'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.913'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.31%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.172'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.52%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06790'
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '71.74'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.49%  '
